{"js": "// The document had three related content edits (per commit message\n// \"Removed limited buffer in the doc\" and the underlying diff):\n//   1. The KPI bullet \"percentage of loss packets (due to overflow in\n//      the queues)\" was reworded to \"packets in buffer over time\".\n//   2. The \"Tx and Rx buffer size\" Factors bullet was removed entirely.\n//   3. The Assumptions bullet \"FIFO queues with limited capacity (=> maybe\n//      M/M/1/C)\" changed \"limited\" to \"unlimited\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// 1) Reword the KPI bullet about loss packets -> packets-in-buffer chart label.\nfor (const p of paragraphs.items) {\n  if (p.text === \"percentage of loss packets (due to overflow in the queues)\") {\n    p.getRange().insertText(\"packets in buffer over time\", Word.InsertLocation.replace);\n    break;\n  }\n}\n\n// 2) Delete the \"Tx and Rx buffer size\" bullet paragraph entirely.\nfor (const p of paragraphs.items) {\n  if (p.text === \"Tx and Rx buffer size\") {\n    p.delete();\n    break;\n  }\n}\n\nawait context.sync();\n\n// 3) Change \"FIFO queues with limited capacity\" -> \"... unlimited capacity\".\nconst matches = body.search(\"FIFO queues with limited capacity\", { matchCase: true });\nmatches.load(\"items\");\nawait context.sync();\n\nfor (const m of matches.items) {\n  m.insertText(\"FIFO queues with unlimited capacity\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The document had three related content edits (per commit message\n# \"Removed limited buffer in the doc\" and the underlying diff):\n#   1. The KPI bullet \"percentage of loss packets (due to overflow in\n#      the queues)\" was reworded to \"packets in buffer over time\".\n#   2. The \"Tx and Rx buffer size\" Factors bullet was removed entirely.\n#   3. The Assumptions bullet \"FIFO queues with limited capacity (=> maybe\n#      M/M/1/C)\" changed \"limited\" to \"unlimited\".\n\n$d = $word.ActiveDocument\n\n# 1) Reword the KPI bullet about loss packets -> packets-in-buffer chart label.\nforeach ($p in $d.Paragraphs) {\n    $txt = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($txt -eq \"percentage of loss packets (due to overflow in the queues)\") {\n        $p.Range.Text = \"packets in buffer over time\"\n        break\n    }\n}\n\n# 2) Delete the \"Tx and Rx buffer size\" bullet paragraph entirely.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $txt = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($txt -eq \"Tx and Rx buffer size\") {\n        $target = $p\n        break\n    }\n}\nif ($target -ne $null) {\n    $target.Range.Delete()\n}\n\n# 3) Change \"FIFO queues with limited capacity\" -> \"... unlimited capacity\".\n$find = $d.Content.Find\n$find.Text = \"FIFO queues with limited capacity\"\n$found = $find.Execute()\nif ($found) {\n    $find.Parent.Text = \"FIFO queues with unlimited capacity\"\n}\n"}
